$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference an untouched cell's style (no explicit number format) so that
# forcing text-entry on numeric-looking "Price" values doesn't leave a
# stray style index behind on the edited cells.
$plainStyle = $ws.Range("D4").Style

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $plainStyle
}

Set-TextValue 'D2' '22.382.06'
$ws.Range('E2').Value = '  -0.71%  '
Set-TextValue 'D3' '1.574.59'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.59%  '
Set-TextValue 'D5' '1.002'
$ws.Range('E5').Value = '  -0.58%  '
Set-TextValue 'D6' '290.62'
$ws.Range('E6').Value = '  -0.73%  '
Set-TextValue 'D7' '0.3758'
$ws.Range('E7').Value = '  +2.56%  '
Set-TextValue 'D8' '50.06'
$ws.Range('E8').Value = '  +1.54%  '
$ws.Range('E9').Value = '  +2.61%  '
$ws.Range('E10').Value = '  +0.54%  '
Set-TextValue 'D11' '0.07682'
$ws.Range('E11').Value = '  +1.50%  '
$ws.Range('E12').Value = '  -0.56%  '
Set-TextValue 'D13' '21.36'
$ws.Range('E13').Value = '  +1.69%  '
Set-TextValue 'D14' '5.991'
$ws.Range('E14').Value = '  -0.03%  '
Set-TextValue 'D15' '6.940'
$ws.Range('E15').Value = '  +1.22%  '
Set-TextValue 'D16' '0.00001144'
$ws.Range('E16').Value = '  +0.83%  '
Set-TextValue 'D17' '1.574.98'
$ws.Range('E17').Value = '  +0.33%  '
Set-TextValue 'D18' '90.40'
$ws.Range('E18').Value = '  +0.92%  '
Set-TextValue 'D19' '0.06724'
$ws.Range('E19').Value = '  -0.93%  '
Set-TextValue 'D20' '1.002'
$ws.Range('E20').Value = '  -0.65%  '
$ws.Range('E21').Value = '  +3.16%  '
Set-TextValue 'D22' '6.245'
$ws.Range('E22').Value = '  +0.52%  '
Set-TextValue 'D23' '0.5272'
$ws.Range('E23').Value = '  -5.05%  '
$ws.Range('E24').Value = '  +1.58%  '
Set-TextValue 'D25' '22.386.98'
$ws.Range('E25').Value = '  -0.67%  '
Set-TextValue 'D26' '2.393'
$ws.Range('E26').Value = '  +0.66%  '
Set-TextValue 'D27' '2.768'
$ws.Range('E27').Value = '  -5.43%  '
Set-TextValue 'D28' '20.25'
$ws.Range('E28').Value = '  +2.71%  '
Set-TextValue 'D29' '144.59'
$ws.Range('E29').Value = '  -0.36%  '
Set-TextValue 'D30' '5.058'
$ws.Range('E30').Value = '  +1.73%  '
Set-TextValue 'D31' '126.48'
$ws.Range('E31').Value = '  +0.99%  '
Set-TextValue 'D32' '1.745.90'
$ws.Range('E32').Value = '  -0.24%  '
Set-TextValue 'D33' '1.028'
$ws.Range('E33').Value = '  +11.09%  '
Set-TextValue 'D34' '6.270'
$ws.Range('E34').Value = '  +0.39%  '
Set-TextValue 'D35' '2.023'
$ws.Range('E35').Value = '  -1.60%  '
Set-TextValue 'D36' '10.14'
$ws.Range('E36').Value = '  -2.60%  '
Set-TextValue 'D37' '0.08561'
$ws.Range('E37').Value = '  +0.08%  '
Set-TextValue 'D38' '0.02548'
$ws.Range('E38').Value = '  +1.85%  '
Set-TextValue 'D39' '0.2329'
$ws.Range('E39').Value = '  +2.39%  '
Set-TextValue 'D40' '0.06535'
$ws.Range('E40').Value = '  +0.86%  '
Set-TextValue 'D41' '5.520'
$ws.Range('E41').Value = '  +2.45%  '
$ws.Range('E42').Value = '  +3.14%  '
Set-TextValue 'D43' '11.67'
$ws.Range('E43').Value = '  +0.01%  '
Set-TextValue 'D44' '0.6436'
$ws.Range('E44').Value = '  +1.97%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-TextValue 'D45' '1.001'
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D46' '14.04'
$ws.Range('E46').Value = '  -2.15%  '
Set-TextValue 'D47' '0.6025'
$ws.Range('E47').Value = '  +2.10%  '
Set-TextValue 'D48' '3.776'
$ws.Range('E48').Value = '  -0.58%  '
Set-TextValue 'D49' '1.309'
$ws.Range('E49').Value = '  +11.30%  '
Set-TextValue 'D50' '2.099'
$ws.Range('E50').Value = '  +0.27%  '
Set-TextValue 'D51' '125.15'
$ws.Range('E51').Value = '  +1.49%  '

Write-Output "Applied cryptos update."
